$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44305, 'Primera', 40, 24000, 24000, 24000, 'Perú', 1200),
    @(3, 44529, 'Primera', 34, 28000, 28000, 28000, 'Perú', 1400),
    @(4, 44356, 'Primera', 100, 20000, 21000, 20500, 'Perú', 1025),
    @(5, 44298, 'Primera', 240, 19000, 20000, 19500, 'Perú', 975),
    @(6, 44372, 'Primera', 60, 20000, 21000, 20667, 'Perú', 1033),
    @(7, 44382, 'Primera', 200, 19000, 20000, 19500, 'Perú', 975),
    @(8, 44620, 'Primera', 60, 22000, 22000, 22000, 'Perú', 1100),
    @(9, 44166, 'Primera', 120, 28000, 28000, 28000, 'Perú', 1400),
    @(10, 44302, 'Primera', 100, 19000, 20000, 19500, 'Perú', 975),
    @(11, 44445, 'Primera', 35, 20000, 20000, 20000, 'Perú', 1000),
    @(12, 44299, 'Primera', 150, 19000, 20000, 19500, 'Perú', 975),
    @(13, 44830, 'Primera', 200, 30000, 30000, 30000, 'Perú', 1500),
    @(14, 44355, 'Primera', 200, 20000, 21000, 20500, 'Ecuador', 1025),
    @(15, 44365, 'Primera', 150, 20000, 21000, 20500, 'Perú', 1025),
    @(16, 44473, 'Primera', 40, 24000, 24000, 24000, 'Perú', 1200),
    @(17, 44350, 'Primera', 90, 21000, 22000, 21556, 'Perú', 1078),
    @(18, 44431, 'Primera', 60, 25000, 25000, 25000, 'Perú', 1250),
    @(19, 44263, 'Segunda', 150, 15000, 15000, 15000, 'Perú', 750),
    @(20, 44165, 'Primera', 300, 27000, 28000, 27500, 'Perú', 1375),
    @(21, 44417, 'Primera', 30, 24000, 24000, 24000, 'Perú', 1200),
    @(22, 44452, 'Primera', 35, 21000, 22000, 21429, 'Perú', 1071),
    @(23, 44333, 'Primera', 30, 22000, 22000, 22000, 'Perú', 1100),
    @(24, 44438, 'Primera', 25, 21000, 21000, 21000, 'Perú', 1050),
    @(25, 44357, 'Primera', 200, 20000, 21000, 20500, 'Perú', 1025),
    @(26, 44760, 'Primera', 300, 24000, 25000, 24500, 'Perú', 1225),
    @(27, 44300, 'Primera', 150, 19000, 20000, 19500, 'Perú', 975),
    @(28, 44354, 'Primera', 150, 21000, 22000, 21500, 'Perú', 1075),
    @(29, 44613, 'Primera', 60, 30000, 30000, 30000, 'Perú', 1500),
    @(30, 44442, 'Primera', 30, 22000, 22000, 22000, 'Perú', 1100),
    @(31, 44326, 'Primera', 40, 22000, 22000, 22000, 'Perú', 1100),
    @(32, 44284, 'Primera', 40, 23000, 23000, 23000, 'Perú', 1150),
    @(33, 45243, 'Primera', 52, 22000, 22000, 22000, 'Perú', 1100),
    @(34, 44424, 'Primera', 70, 24000, 25000, 24429, 'Perú', 1221),
    @(35, 44410, 'Primera', 40, 25000, 25000, 25000, 'Perú', 1250),
    @(36, 44522, 'Primera', 25, 30000, 30000, 30000, 'Perú', 1500),
    @(37, 44312, 'Primera', 50, 22000, 22000, 22000, 'Perú', 1100),
    @(38, 44270, 'Primera', 50, 24000, 24000, 24000, 'Perú', 1200),
    @(39, 44277, 'Primera', 60, 24000, 24000, 24000, 'Perú', 1200),
    @(40, 44396, 'Primera', 45, 22000, 22000, 22000, 'Perú', 1100),
    @(41, 44363, 'Primera', 150, 21000, 22000, 21500, 'Perú', 1075),
    @(42, 44435, 'Primera', 60, 25000, 25000, 25000, 'Perú', 1250)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 12).Value = $row[2]
    $ws.Cells.Item($r, 13).Value = $row[3]
    $ws.Cells.Item($r, 14).Value = $row[4]
    $ws.Cells.Item($r, 15).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 18).Value = $row[7]
    $ws.Cells.Item($r, 19).Value = $row[8]
}

Write-Output "done"